$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 8, shifting old row 8.. down to row 11..
$ws.Rows("8:10").Insert()

# Write column-A labels first, then column-B bodies, so new shared-string
# entries land in the same append order as the target workbook.
$ws.Cells.Item(8, 1).Value = "e006a"
$ws.Cells.Item(9, 1).Value = "e006b"
$ws.Cells.Item(10, 1).Value = "e006b"

$ws.Cells.Item(8, 2).Value = @"
<Bold>e006a Retrofit Period</Bold> 
<InlineUIContainer><Button Content='r27.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
During the periods on the Combat <InlineUIContainer><Button Content='Calendar' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
 marked Refitting, the Division is building itself back up to strength, replacing tanks, and retraining new crews. 
<LineBreak/><LineBreak/>
During a refit period, you have the option of replacing your current tank per 
<InlineUIContainer><Button Content='r24.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. 
<LineBreak/><LineBreak/>

"@

$ws.Cells.Item(9, 2).Value = @"
<Bold>e006b Retrofit Period - Crew Training</Bold> 
<InlineUIContainer><Button Content='r27.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<InlineUIContainer><Button Content='r27.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Since the refit takes at least 7 days, you may attempt to improve your crew ratings per 
<InlineUIContainer><Button Content='r27.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
and train your crew to use the gyrostablilizer per 
<InlineUIContainer><Button Content='r27.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
"@

$ws.Cells.Item(10, 2).Value = @"
<Bold>e006b Retrofit Period - Gyrostabilizer</Bold> 
<InlineUIContainer><Button Content='r27.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
When the refit takes at least 7 days, and your crew has a combined rating of 30, your crew is trained in the use of the gyrostablilizer per 
<InlineUIContainer><Button Content='r27.21' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. Loss of the gunner losses the Horizontal Volute Spring Suspension (HVSS) capability. Click image to continue.
<LineBreak/><LineBreak/>
                                                  <InlineUIContainer><Image Name='c75Hvss'  Height='80' Width='80'></Image></InlineUIContainer>
"@

$ws.Rows.Item(8).RowHeight = 150
$ws.Rows.Item(9).RowHeight = 105
$ws.Rows.Item(10).RowHeight = 105

# Update the sheet view to match the author's final cursor position
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("B11:B12").Select()
$ws.Application.ActiveCell.Activate()
